$d = $word.ActiveDocument
$d.Content.Find.Execute("43×84=3612", $true, $false, $false, $false, $false, $true, 1, $false, "35×36=1260", 2)
$d.Content.Find.Execute("17×82=1394", $true, $false, $false, $false, $false, $true, 1, $false, "76×92=6992", 2)
$d.Content.Find.Execute("94×83=7802", $true, $false, $false, $false, $false, $true, 1, $false, "23×86=1978", 2)
$d.Content.Find.Execute("96×81=7776", $true, $false, $false, $false, $false, $true, 1, $false, "91×68=6188", 2)
$d.Content.Find.Execute("77×82=6314", $true, $false, $false, $false, $false, $true, 1, $false, "80×67=5360", 2)
$d.Content.Find.Execute("62×17=1054", $true, $false, $false, $false, $false, $true, 1, $false, "62×51=3162", 2)
$d.Content.Find.Execute("23×17=391", $true, $false, $false, $false, $false, $true, 1, $false, "19×97=1843", 2)
$d.Content.Find.Execute("81×49=3969", $true, $false, $false, $false, $false, $true, 1, $false, "96×89=8544", 2)
$d.Content.Find.Execute("33×42=1386", $true, $false, $false, $false, $false, $true, 1, $false, "56×70=3920", 2)
$d.Content.Find.Execute("22×86=1892", $true, $false, $false, $false, $false, $true, 1, $false, "12×58=696", 2)
$d.Content.Find.Execute("13×23=299", $true, $false, $false, $false, $false, $true, 1, $false, "99×34=3366", 2)
$d.Content.Find.Execute("49×35=1715", $true, $false, $false, $false, $false, $true, 1, $false, "26×50=1300", 2)
$d.Content.Find.Execute("21×25=525", $true, $false, $false, $false, $false, $true, 1, $false, "56×99=5544", 2)
$d.Content.Find.Execute("78×43=3354", $true, $false, $false, $false, $false, $true, 1, $false, "83×33=2739", 2)
$d.Content.Find.Execute("72×63=4536", $true, $false, $false, $false, $false, $true, 1, $false, "43×69=2967", 2)
$d.Content.Find.Execute("80×60=4800", $true, $false, $false, $false, $false, $true, 1, $false, "60×71=4260", 2)
$d.Content.Find.Execute("35×27=945", $true, $false, $false, $false, $false, $true, 1, $false, "49×29=1421", 2)
$d.Content.Find.Execute("29×46=1334", $true, $false, $false, $false, $false, $true, 1, $false, "75×48=3600", 2)
$d.Content.Find.Execute("32×55=1760", $true, $false, $false, $false, $false, $true, 1, $false, "61×26=1586", 2)
$d.Content.Find.Execute("32×28=896", $true, $false, $false, $false, $false, $true, 1, $false, "11×53=583", 2)
$d.Content.Find.Execute("70×45=3150", $true, $false, $false, $false, $false, $true, 1, $false, "48×31=1488", 2)
$d.Content.Find.Execute("20×88=1760", $true, $false, $false, $false, $false, $true, 1, $false, "29×73=2117", 2)
$d.Content.Find.Execute("83×90=7470", $true, $false, $false, $false, $false, $true, 1, $false, "54×46=2484", 2)
$d.Content.Find.Execute("62×38=2356", $true, $false, $false, $false, $false, $true, 1, $false, "63×20=1260", 2)
$d.Content.Find.Execute("86×32=2752", $true, $false, $false, $false, $false, $true, 1, $false, "13×25=325", 2)
